$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 11 (pushes nothing down, just extends the table by one
# row) - inheriting the number formatting from the row above it (row 10),
# exactly like typing a new entry under the existing date/hours table.
$ws.Rows("11").Insert(-4121)

# New appropriation entry: date 2013-10-16, 2h15min (0.09375 of a day).
$ws.Range("A11").Value = 41563
$ws.Range("B11").Value = 0.09375

# Move the active selection to C11, matching where the user would land
# after filling in the new row.
$ws.Range("C11").Select()
